$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test cases")

# --- Row 29: new "Function 07" section header -----------------------------
# Copy the formatting from row 11 (an existing section-header row that uses
# the exact same style pattern: s4,s4,s4,s7,s7,s7,s7,s7,s8,s4) onto row 29,
# then set the header text.
$ws.Range("A11:J11").Copy()
$ws.Range("A29:J29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Function 07: Progress Tracking"

# --- Row 30: S-PT-01 / UI01 test case --------------------------------------
# Row 5 has the matching style pattern (s1, default, s2 x6, s3) and the same
# row height (60), so copy its formats first.
$ws.Range("A5:I5").Copy()
$ws.Range("A30:I30").PasteSpecial(-4122)

$ws.Range("A30").Value = "S-PT-01"
$ws.Range("B30").Value = "UI01"
$ws.Range("C30").Value = "User click ""See achievements"" button, cahnge layout to Achivement View"
$ws.Range("D30").Value = "1. Currently in DiaryScene" + [char]10 + "2. Click ""See achievements"" button"
$ws.Range("E30").Value = "The layout change to Achievement view"
$ws.Range("F30").Value = "The same ER"
$ws.Range("G30").Value = "pass"
$ws.Range("H30").Value = "LDBach"
$ws.Range("I30").Value = "17/12/207"

# --- Row 31: S-PT-01 / UI02 test case --------------------------------------
$ws.Range("A31:I31").Copy()
$ws.Range("A31:I31").PasteSpecial(-4122)

$ws.Range("A31").Value = "S-PT-01"
$ws.Range("B31").Value = "UI02"
$ws.Range("C31").Value = "User click ""See pet's desire"" button, cahnge layout to Achivement View"
$ws.Range("D31").Value = "1. Currently in DiaryScene" + [char]10 + "2. Click ""See achievements"" button" + [char]10 + "3. Wait for Achievement screen come up" + [char]10 + "4. Click ""See Pet's Desire"" button"
$ws.Range("E31").Value = "The layout change to Achievement view, and turn back to Quest View after clicking ""See Pet's Desire"" button"
$ws.Range("F31").Value = "The same ER"
$ws.Range("G31").Value = "Pass"
$ws.Range("H31").Value = "LDBach"
$ws.Range("I31").Value = "17/12/2017"

# --- sheet view: scroll position + selection -------------------------------
$ws.Range("I32").Select()
$excel.ActiveWindow.ScrollRow = 26
